# Update crypto price/volume data per the scraped snapshot refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.134.70'
$ws.Range('E2').Value = '  -0.57%  '
$ws.Range('D3').Value = '2.027.16'
$ws.Range('E3').Value = '  -1.03%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = "'226.70"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.05%  '
$ws.Range('D6').Value = "'0.610"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.80%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').Value = "'55.09"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.79%  '
$ws.Range('E9').Value = '  -1.90%  '
$ws.Range('D10').Value = "'0.0785"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.20%  '
$ws.Range('E11').Value = '  -4.78%  '
$ws.Range('D12').Value = '2.315.39'
$ws.Range('E12').Value = '  -0.64%  '
$ws.Range('E13').Value = '  -3.45%  '
$ws.Range('D14').Value = "'20.22"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.04%  '
$ws.Range('D15').Value = "'0.743"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.60%  '
$ws.Range('E16').Value = '  -2.26%  '
$ws.Range('D17').Value = '2.022.09'
$ws.Range('E17').Value = '  -1.07%  '
$ws.Range('D18').Value = '37.080.36'
$ws.Range('E18').Value = '  -0.42%  '
$ws.Range('D19').Value = "'6.51"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +6.96%  '
$ws.Range('D20').Value = "'68.80"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.86%  '
$ws.Range('D21').Value = '0.0₃0816'
$ws.Range('E21').Value = '  -1.20%  '
$ws.Range('D22').Value = "'222.94"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.28%  '
$ws.Range('E23').Value = '  +0.06%  '
$ws.Range('D24').Value = "'2.45"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.05%  '
$ws.Range('D25').Value = "'2.19"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -4.12%  '
$ws.Range('D26').Value = "'165.64"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.93%  '
$ws.Range('D27').Value = "'9.21"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -5.29%  '
$ws.Range('E28').Value = '  -0.52%  '
$ws.Range('D29').Value = "'18.68"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E30').Value = '  -2.42%  '
$ws.Range('D31').Value = "'0.118"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.23%  '
$ws.Range('D32').Value = "'4.51"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.52%  '
$ws.Range('E33').Value = '  -0.94%  '
$ws.Range('D34').Value = "'4.47"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.30%  '
$ws.Range('E35').Value = '  -2.94%  '
$ws.Range('D36').Value = "'1.86"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.79%  '
$ws.Range('E37').Value = '  +0.00%  '
$ws.Range('B38').Value = 'THORChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D38').Value = "'5.54"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +5.52%  '
$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D39').Value = "'3.11"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -3.85%  '
$ws.Range('D40').Value = '1.469.11'
$ws.Range('E40').Value = '  -0.47%  '
$ws.Range('E41').Value = '  -2.99%  '
$ws.Range('D42').Value = "'95.59"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.83%  '
$ws.Range('E43').Value = '  -3.23%  '
$ws.Range('D44').Value = "'16.42"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -3.35%  '
$ws.Range('D45').Value = "'0.0907"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -3.85%  '
$ws.Range('E46').Value = '  -2.33%  '
$ws.Range('D47').Value = "'7.27"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.34%  '
$ws.Range('E48').Value = '  -0.90%  '
$ws.Range('E49').Value = '  +0.24%  '
$ws.Range('B50').Value = 'RocketPoolETH'
$ws.Range('C50').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D50').Value = '2.208.63'
$ws.Range('E50').Value = '  -1.07%  '
$ws.Range('B51').Value = 'FTXToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D51').Value = "'3.64"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -7.38%  '
